$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# Helper: wrap a WordprocessingML body fragment in the minimal package
# envelope that Range.InsertXML expects, and insert it into a Range,
# fully replacing the range content without leaving stray formatting-only
# elements (such as proofErr marks) behind.
# NOTE: always pre-build the xml string into its own variable before
# calling this (avoid parenthesised expressions as call arguments next to
# a COM object argument - this runtime mis-parses that combination).
# -------------------------------------------------------------------------
function Insert-BodyXml($range, $bodyXml) {
    $pkgHead = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
    $pkgTail = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $pkg = $pkgHead + $bodyXml + $pkgTail
    $range.InsertXML($pkg)
}

# -------------------------------------------------------------------------
# 1) "		./audio_rec | ./audio_play" paragraph: drop the spell/grammar proof
#    marks and collapse the runs into a single simple run (same visible
#    text), keeping the two leading tabs split exactly as in the target.
# -------------------------------------------------------------------------
$pAudio = $d.Paragraphs(5)
$pAudio.Range.Delete()
$pAudio = $d.Paragraphs(5)
$pAudio.Range.InsertParagraphBefore()
$pAudio = $d.Paragraphs(5)
$audioXml = '<w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>./audio_rec | ./audio_play</w:t></w:r></w:p>'
$audioRange = $pAudio.Range
Insert-BodyXml $audioRange $audioXml

# -------------------------------------------------------------------------
# 2) "There was a considerable delay." -> "There was a little delay.",
#    with the (moved) _GoBack bookmark now sitting right after "little".
#    First drop the pre-existing _GoBack bookmark (it currently lives a
#    few paragraphs further down) so we don't end up with two of them.
# -------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$pDelay = $d.Paragraphs(8)
$pDelay.Range.Delete()
$pDelay = $d.Paragraphs(8)
$pDelay.Range.InsertParagraphBefore()
$pDelay = $d.Paragraphs(8)
$delayXml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">There was a </w:t></w:r><w:r><w:t>little</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> delay.</w:t></w:r></w:p>'
$delayRange = $pDelay.Range
Insert-BodyXml $delayRange $delayXml

# -------------------------------------------------------------------------
# 3) Replace the bookmark-only numbered paragraph + blank paragraph +
#    four "The problems you found / ..." bullet paragraphs (6 paragraphs)
#    with four new indented dash paragraphs describing the WiFi issues.
# -------------------------------------------------------------------------
$pStart = $d.Paragraphs(11)
$pEnd = $d.Paragraphs(16)
$fullRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$fullRange.Delete()

$pTarget = $d.Paragraphs(11)
$wifiXml = '<w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>- When connected to WiFi</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:tab/><w:t>-  Pulse break</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:tab/><w:t>-  Huge delays</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:tab/><w:t>-  More noise</w:t></w:r></w:p>'
$targetRange = $pTarget.Range
Insert-BodyXml $targetRange $wifiXml
